$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Invoice")

# E18 previously held the numeric value 1 (discount amount placeholder).
# It is replaced with explanatory text, which makes the dependent formula
# in F18 (and the totals that roll up from it) evaluate to #VALUE!.
$ws.Range("E18").Value = "This client doesn't benefit from any discount"

# Footer note text updated.
$ws.Range("A31").Value = "Bla Bla Bla"
